# phase1-A: data counts frames
# The sheet's column A header (A1) was blank; label it "party" so the
# row headers (polling stations / E-votes / Total) are identified.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("A1").Value = "party"

# Leave the cursor where the author left it after the edit.
$ws.Range("B10").Select()
